$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells, copying the existing header style (AC1)
# over to AD1:AF1 first so they pick up the same bold/bordered header
# formatting used by the rest of row 1, then set their text.
$ws.Range("AC1:AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value = 74
    $ws.Cells.Item($r, 31).Value = 88
    $ws.Cells.Item($r, 32).Value = 0
}
